# "task creation bug resolve"
#
# The whole task-creation flow diagram on slide 2 (all rectangles, the
# "TASK CREATION" box, the elbow connectors and the floating caption text
# boxes) was nudged by a fixed amount: +279400 EMU on the x-axis and
# -8467 EMU on the y-axis (+22pt right, ~0.667pt up). Apply that same
# delta to every shape on the slide so the diagram keeps its exact
# internal layout/connections, just shifted as a whole.

$EMU_PER_POINT = 12700
$dxEmu = 279400
$dyEmu = -8467

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)

    # PowerPoint's COM surface works in points; round-trip through EMU
    # (the unit the OOXML <a:off> actually stores) so the shift lands on
    # the exact target EMU value instead of drifting from float error.
    $curLeftEmu = [Math]::Round($shp.Left * $EMU_PER_POINT)
    $curTopEmu = [Math]::Round($shp.Top * $EMU_PER_POINT)

    $newLeftEmu = $curLeftEmu + $dxEmu
    $newTopEmu = $curTopEmu + $dyEmu

    # Nudge by half an EMU before converting back to points so the
    # engine's internal float truncation still resolves to the exact
    # integer EMU we want.
    $shp.Left = ($newLeftEmu + 0.5) / $EMU_PER_POINT
    $shp.Top = ($newTopEmu + 0.5) / $EMU_PER_POINT
}
